$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: rotate the c3/c4/c5 placeholders into real a5/b5/c5 values ---
# Set in this order (C5, B5, A5) so new shared strings are appended in the
# same order the target workbook uses (c5 reuses the existing string,
# b5 becomes the next new unique string, a5 the one after).
$ws.Range("C5").Value = "c5"
$ws.Range("B5").Value = "b5"
$ws.Range("A5").Value = "a5"

# --- Row 3: new numeric cell D3 and wrapped text cell F3 ---
$ws.Range("D3").Value = 123456

$ws.Range("F3").Value = "`nwhitespace    string"
$ws.Range("F3").WrapText = $true

# Row 3 gets taller to show the wrapped text, column F gets wider.
$ws.Rows.Item(3).RowHeight = 29.25
$ws.Columns.Item(6).ColumnWidth = 27.7109375

# --- Selection moves to C4 ---
$null = $ws.Range("C4").Select()
